$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "TABLE" mini-schema is being renamed to "seat", and a new
# "seat_total_cost" attribute row is inserted right below the existing
# seat/table attribute rows (pushing the PEOPLE block down by one row).

# 1) Insert a new row at row 8 - this shifts PEOPLE..people_point (rows 9-13)
#    down to rows 10-14 and keeps their formatting intact.
$ws.Rows("8:8").Insert()

# 2) Rename the TABLE entity/attributes to the seat equivalents (in the
#    same order the workbook's shared-string table was authored in).
$ws.Range("D2").Value = "seat"
$ws.Range("D4").Value = "seat_capacity"
$ws.Range("D5").Value = "seat_people_num"
$ws.Range("D6").Value = "seat_time_first_order"
$ws.Range("D7").Value = "seat_master (default table_num)"
$ws.Range("F4").Value = "seat_num(fk)"
$ws.Range("D3").Value = "seat_id(pk)"

# 3) Populate the newly inserted row with the new attribute.
$ws.Range("D8").Value = "seat_total_cost "

# 4) Update the active selection to match the authored workbook state.
$ws.Range("D8").Select()
